$wb = $excel.ActiveWorkbook

# --- Finance sheet: add row 3 with date + full medians row ---
$wsFinance = $wb.Worksheets.Item("Finance")
$wsFinance.Range("A3").NumberFormat = "@"
$wsFinance.Range("A3").Value = "2025-10-22"
$wsFinance.Range("A3").ClearFormats()
$wsFinance.Range("B3").Value = 5.68
$wsFinance.Range("C3").Value = 9.32
$wsFinance.Range("D3").Value = 7.66
$wsFinance.Range("E3").Value = 5.03
$wsFinance.Range("F3").Value = 26.84
$wsFinance.Range("G3").Value = 2.73
$wsFinance.Range("H3").Value = 5.36
$wsFinance.Range("I3").Value = 3.5
$wsFinance.Range("J3").Value = 2.73
$wsFinance.Range("K3").Value = 12.58
$wsFinance.Range("L3").Value = 5.83
$wsFinance.Range("M3").Value = 5.68

# --- Non-Finance sheet: add row 3 with just the date ---
$wsNonFinance = $wb.Worksheets.Item("Non-Finance")
$wsNonFinance.Range("A3").NumberFormat = "@"
$wsNonFinance.Range("A3").Value = "2025-10-22"
$wsNonFinance.Range("A3").ClearFormats()
